# Apply settings changes per commit:
# - switch all testing to parametric
# - add missing collector$p.adj.signif to parametric testing

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# project_name
$ws.Range("B3").Value = "dc10"

# excluded_channels
$ws.Range("B5").Value = "B2M, DNA, Bead, LD, Live, Dead, ID, Cell-ID, Cell_ID, NA"

# anchor_ids
$ws.Range("B9").Value = "Anchor"

# grouping_columns
$ws.Range("B15").Value = "group, paired_0, paired_0_5, paired_4, paired_24, paired_0_LPS, paired_0_5_LPS, paired_4_LPS, paired_24_LPS"

# grouping_orders
$ws.Range("B16").Value = "DC_0h, DC_0_5h, DC_4h, DC_24h, DC_LPS_0h, DC_LPS_0_5h, DC_LPS_4h, DC_LPS_24h, DC10_0h, DC10_0_5h, DC10_4h, DC10_24h, DC10_LPS_0h, DC10_LPS_0_5h, DC10_LPS_4h, DC10_LPS_24h; DC_0h, DC10_0h; DC_0_5h, DC10_0_5h; DC_4h, DC10_4h; DC_24h, DC10_24h; DC_LPS_0h, DC10_LPS_0h; DC_LPS_0_5h, DC10_LPS_0_5h; DC_LPS_4h, DC10_LPS_4h; DC_LPS_24h, DC10_LPS_24h"

# data_subsets
$ws.Range("B18").Value = "DC"

# event_cutoff
$ws.Range("B28").Value = 20

# ccp_delta_cutoff
$ws.Range("B30").Value = 0.025

# umap_n
$ws.Range("B32").Value = 20

# umap_min_dist
$ws.Range("B33").Value = 0.15

# Row 16 grows taller to fit the longer wrapped grouping_orders text
$ws.Rows.Item(16).RowHeight = 86.4

# Update view to match: scrolled up, selection on B16
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("B16").Select()
